$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.646881333333333
$ws.Range("H2").Value = 4.940644
$ws.Range("I2").Value = 0.3367300927127475
$ws.Range("J2").Value = 0.3367300927127475
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1503136666666667
$ws.Range("N2").Value = 0.450941
$ws.Range("O2").Value = 0.008029526741163598
$ws.Range("P2").Value = 0.008029526741163598
$ws.Range("Q2").Value = 0.2475487717782222
$ws.Range("R2").Value = 2.227938946004
$ws.Range("S2").Value = 0.002703783283991504
$ws.Range("T2").Value = 0.002703783283991504

# Row 3
$ws.Range("G3").Value = 1.646881333333333
$ws.Range("H3").Value = 4.940644
$ws.Range("I3").Value = 0.3367300927127475
$ws.Range("J3").Value = 0.3367300927127475
$ws.Range("O3").Value = 0.8389317081486641
$ws.Range("P3").Value = 0.8389317081486641
$ws.Range("Q3").Value = 25.86410390706444
$ws.Range("R3").Value = 232.77693516358
$ws.Range("S3").Value = 0.2824935518645633
$ws.Range("T3").Value = 0.2824935518645633

# Row 4
$ws.Range("G4").Value = 1.646881333333333
$ws.Range("H4").Value = 4.940644
$ws.Range("I4").Value = 0.3367300927127475
$ws.Range("J4").Value = 0.3367300927127475
$ws.Range("M4").Value = 2.758466666666667
$ws.Range("N4").Value = 8.2754
$ws.Range("O4").Value = 0.1473530807662759
$ws.Range("P4").Value = 0.1473530807662759
$ws.Range("Q4").Value = 4.542867261955555
$ws.Range("R4").Value = 40.8858053576
$ws.Range("S4").Value = 0.04961821654793706
$ws.Range("T4").Value = 0.04961821654793706

# Row 5
$ws.Range("G5").Value = 1.646881333333333
$ws.Range("H5").Value = 4.940644
$ws.Range("I5").Value = 0.3367300927127475
$ws.Range("J5").Value = 0.3367300927127475
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1064366666666667
$ws.Range("N5").Value = 0.31931
$ws.Range("O5").Value = 0.005685684343896315
$ws.Range("P5").Value = 0.005685684343896314
$ws.Range("Q5").Value = 0.1752885595155555
$ws.Range("R5").Value = 1.57759703564
$ws.Range("S5").Value = 0.001914541016255623
$ws.Range("T5").Value = 0.001914541016255623

# Row 6
$ws.Range("I6").Value = 0.3135726931406526
$ws.Range("J6").Value = 0.3135726931406525
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1503136666666667
$ws.Range("N6").Value = 0.450941
$ws.Range("O6").Value = 0.008029526741163598
$ws.Range("P6").Value = 0.008029526741163598
$ws.Range("Q6").Value = 0.2305244964143334
$ws.Range("R6").Value = 2.074720467729
$ws.Range("S6").Value = 0.002517840324871557
$ws.Range("T6").Value = 0.002517840324871556

# Row 7
$ws.Range("I7").Value = 0.3135726931406526
$ws.Range("J7").Value = 0.3135726931406525
$ws.Range("O7").Value = 0.8389317081486641
$ws.Range("P7").Value = 0.8389317081486641
$ws.Range("S7").Value = 0.2630660750852645
$ws.Range("T7").Value = 0.2630660750852645

# Row 8
$ws.Range("I8").Value = 0.3135726931406526
$ws.Range("J8").Value = 0.3135726931406525
$ws.Range("M8").Value = 2.758466666666667
$ws.Range("N8").Value = 8.2754
$ws.Range("O8").Value = 0.1473530807662759
$ws.Range("P8").Value = 0.1473530807662759
$ws.Range("Q8").Value = 4.230447924733334
$ws.Range("R8").Value = 38.0740313226
$ws.Range("S8").Value = 0.04620590237845323
$ws.Range("T8").Value = 0.04620590237845322

# Row 9
$ws.Range("I9").Value = 0.3135726931406526
$ws.Range("J9").Value = 0.3135726931406525
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1064366666666667
$ws.Range("N9").Value = 0.31931
$ws.Range("O9").Value = 0.005685684343896315
$ws.Range("P9").Value = 0.005685684343896314
$ws.Range("Q9").Value = 0.1632337200433334
$ws.Range("R9").Value = 1.46910348039
$ws.Range("S9").Value = 0.001782875352063212
$ws.Range("T9").Value = 0.001782875352063211

# Row 10
$ws.Range("G10").Value = 1.115861333333333
$ws.Range("H10").Value = 3.347584
$ws.Range("I10").Value = 0.2281549269050169
$ws.Range("J10").Value = 0.2281549269050168
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1503136666666667
$ws.Range("N10").Value = 0.450941
$ws.Range("O10").Value = 0.008029526741163598
$ws.Range("P10").Value = 0.008029526741163598
$ws.Range("Q10").Value = 0.1677292085048889
$ws.Range("R10").Value = 1.509562876544
$ws.Range("S10").Value = 0.001831976086712059
$ws.Range("T10").Value = 0.001831976086712059

# Row 11
$ws.Range("G11").Value = 1.115861333333333
$ws.Range("H11").Value = 3.347584
$ws.Range("I11").Value = 0.2281549269050169
$ws.Range("J11").Value = 0.2281549269050168
$ws.Range("O11").Value = 0.8389317081486641
$ws.Range("P11").Value = 0.8389317081486641
$ws.Range("Q11").Value = 17.52448879409778
$ws.Range("R11").Value = 157.72039914688
$ws.Range("S11").Value = 0.1914064025509594
$ws.Range("T11").Value = 0.1914064025509594

# Row 12
$ws.Range("G12").Value = 1.115861333333333
$ws.Range("H12").Value = 3.347584
$ws.Range("I12").Value = 0.2281549269050169
$ws.Range("J12").Value = 0.2281549269050168
$ws.Range("M12").Value = 2.758466666666667
$ws.Range("N12").Value = 8.2754
$ws.Range("O12").Value = 0.1473530807662759
$ws.Range("P12").Value = 0.1473530807662759
$ws.Range("Q12").Value = 3.078066292622222
$ws.Range("R12").Value = 27.7025966336
$ws.Range("S12").Value = 0.03361933137145873
$ws.Range("T12").Value = 0.03361933137145872

# Row 13
$ws.Range("G13").Value = 1.115861333333333
$ws.Range("H13").Value = 3.347584
$ws.Range("I13").Value = 0.2281549269050169
$ws.Range("J13").Value = 0.2281549269050168
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1064366666666667
$ws.Range("N13").Value = 0.31931
$ws.Range("O13").Value = 0.005685684343896315
$ws.Range("P13").Value = 0.005685684343896314
$ws.Range("Q13").Value = 0.1187685607822222
$ws.Range("R13").Value = 1.06891704704
$ws.Range("S13").Value = 0.001297216895886663
$ws.Range("T13").Value = 0.001297216895886662

# Row 14
$ws.Range("G14").Value = 0.5944396666666667
$ws.Range("H14").Value = 1.783319
$ws.Range("I14").Value = 0.1215422872415831
$ws.Range("J14").Value = 0.1215422872415831
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.1503136666666667
$ws.Range("N14").Value = 0.450941
$ws.Range("O14").Value = 0.008029526741163598
$ws.Range("P14").Value = 0.008029526741163598
$ws.Range("Q14").Value = 0.08935240590877778
$ws.Range("R14").Value = 0.8041716531790001
$ws.Range("S14").Value = 0.000975927045588479
$ws.Range("T14").Value = 0.0009759270455884789

# Row 15
$ws.Range("G15").Value = 0.5944396666666667
$ws.Range("H15").Value = 1.783319
$ws.Range("I15").Value = 0.1215422872415831
$ws.Range("J15").Value = 0.1215422872415831
$ws.Range("O15").Value = 0.8389317081486641
$ws.Range("P15").Value = 0.8389317081486641
$ws.Range("Q15").Value = 9.335614530300555
$ws.Range("R15").Value = 84.020530772705
$ws.Range("S15").Value = 0.1019656786478769
$ws.Range("T15").Value = 0.1019656786478769

# Row 16
$ws.Range("G16").Value = 0.5944396666666667
$ws.Range("H16").Value = 1.783319
$ws.Range("I16").Value = 0.1215422872415831
$ws.Range("J16").Value = 0.1215422872415831
$ws.Range("M16").Value = 2.758466666666667
$ws.Range("N16").Value = 8.2754
$ws.Range("O16").Value = 0.1473530807662759
$ws.Range("P16").Value = 0.1473530807662759
$ws.Range("Q16").Value = 1.639742005844445
$ws.Range("R16").Value = 14.7576780526
$ws.Range("S16").Value = 0.0179096304684269
$ws.Range("T16").Value = 0.0179096304684269

# Row 17
$ws.Range("G17").Value = 0.5944396666666667
$ws.Range("H17").Value = 1.783319
$ws.Range("I17").Value = 0.1215422872415831
$ws.Range("J17").Value = 0.1215422872415831
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1064366666666667
$ws.Range("N17").Value = 0.31931
$ws.Range("O17").Value = 0.005685684343896315
$ws.Range("P17").Value = 0.005685684343896314
$ws.Range("Q17").Value = 0.06327017665444444
$ws.Range("R17").Value = 0.56943158989
$ws.Range("S17").Value = 0.0006910510796908181
$ws.Range("T17").Value = 0.0006910510796908179
